$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = '"parent07, "parent05'
$ws.Range("A7").Select()
